# someText.txt was modified: the original "SomeText" paragraph gets
# spell-check proofing marks around it, and a new paragraph "NewText" is
# appended after it (inheriting the trailing "_GoBack" bookmark that used
# to sit at the end of the first paragraph).

$d = $word.ActiveDocument

# --- Step 1: wrap the existing "SomeText" run in <w:proofErr> spellStart/
# spellEnd markers, as Word's background spell checker would after an edit.
$someTextRange = $d.Range(0, 8)
$wrappedRunXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>SomeText</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$someTextRange.InsertXML($wrappedRunXml)

# --- Step 2: the "_GoBack" bookmark currently still wraps "SomeText" in the
# first paragraph; remove it from there so it can move to the new paragraph
# (Word re-drops "_GoBack" at the location of the most recent edit).
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# --- Step 3: append a new paragraph with "NewText", carrying the
# relocated "_GoBack" bookmark at its end, matching paragraph 1's formatting.
$endOfDoc = $d.Content
$endOfDoc.Collapse(0)
$newParaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>NewText</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$endOfDoc.InsertXML($newParaXml)
